$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B14 was stored as text "1"; convert it to a real number 1 (matches the
# rest of the politeness_score column, which is numeric).
$ws.Cells.Item(14, 2).Value = 1

# Append new row 15 with annotation data.
$ws.Cells.Item(15, 1).Value = "Ying Tang"

# politeness_score for this row is stored as text "4" (not a number), so
# force a text number format before assigning to avoid Excel's automatic
# numeric coercion, then drop back to the default style so no stray
# formatting is left on the cell.
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "4"
$ws.Cells.Item(15, 2).Style = "Normal"

$ws.Cells.Item(15, 3).Value = "I suggest"
$ws.Cells.Item(15, 4).Value = "SUG"
$ws.Cells.Item(15, 5).Value = "WRI"
$ws.Cells.Item(15, 6).Value = "3a6bf25f-9f71-48b7-a40b-7e968e5f9337"
$ws.Cells.Item(15, 7).Value = "ry-TW-WAb_annotated.xlsx"
$ws.Cells.Item(15, 8).Value = "I suggest to change it to e.g. 'from the true to the approximate posterior' to avoid confusion."
